# Add a new "demo_hier" dictionary column (D) with hierarchy-coded labels
# that mirror the existing "demo_code" (B) labels, e.g. "Age" -> "1.0. Age",
# "18-29 years" -> "1.1.0. 18-29 years", etc.
#
# The values are written in the same order they were originally authored
# (all of the "leaf" rows first, top to bottom, followed by the section
# "N.0. <Category>" header rows, top to bottom) so that the resulting
# shared-string table is built up in that exact sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header + "leaf" rows (in row order) ---
$ws.Cells.Item(1, 4).Value = 'demo_hier'
$ws.Cells.Item(2, 4).Value = '0. TOTAL DEMOGRAPHICS'
$ws.Cells.Item(4, 4).Value = '1.1.0. 18-29 years'
$ws.Cells.Item(5, 4).Value = '1.2.0. 30-39 years'
$ws.Cells.Item(6, 4).Value = '1.3.0. 40-49 years'
$ws.Cells.Item(7, 4).Value = '1.4.0. 50-59 years'
$ws.Cells.Item(8, 4).Value = '1.5.0. 60+ years'
$ws.Cells.Item(10, 4).Value = '2.1.0. 1'
$ws.Cells.Item(11, 4).Value = '2.2.0. 2'
$ws.Cells.Item(12, 4).Value = '2.3.0. 3'
$ws.Cells.Item(13, 4).Value = '2.4.0. 4+'
$ws.Cells.Item(15, 4).Value = '3.1.0. Children <6 years'
$ws.Cells.Item(16, 4).Value = '3.2.0. Children 6-14 years'
$ws.Cells.Item(17, 4).Value = '3.3.0. Children 15-19 years'
$ws.Cells.Item(18, 4).Value = '3.4.0. No Children <20 years (Rest)'
$ws.Cells.Item(20, 4).Value = '4.1.0. Employed full-time'
$ws.Cells.Item(21, 4).Value = '4.2.0. Employed part-time'
$ws.Cells.Item(22, 4).Value = '4.3.0. Not employed'
$ws.Cells.Item(24, 4).Value = '5.1.0. Low Affluency'
$ws.Cells.Item(25, 4).Value = '5.2.0. Below average affluency'
$ws.Cells.Item(26, 4).Value = '5.3.0. Above average affluency'
$ws.Cells.Item(27, 4).Value = '5.4.0. High affluency'
$ws.Cells.Item(29, 4).Value = '6.1.0. Young singles/couples without'
$ws.Cells.Item(30, 4).Value = '6.2.0. Mature singles/couples without'
$ws.Cells.Item(31, 4).Value = '6.3.0. Senior singles/couples without'
$ws.Cells.Item(32, 4).Value = '6.4.0. Retired singles/couples withou'
$ws.Cells.Item(33, 4).Value = '6.5.0. Family with children <18'
$ws.Cells.Item(34, 4).Value = '6.5.1.0. Family with children <6 years'
$ws.Cells.Item(35, 4).Value = '6.5.1.1.0. Family with children 0-3 years'
$ws.Cells.Item(36, 4).Value = '6.5.1.2.0. Family with children 4-5 years'
$ws.Cells.Item(37, 4).Value = '6.5.2.0. Family with children 6-12 year'
$ws.Cells.Item(38, 4).Value = '6.5.3.0. Family with children 13-18 yea'
$ws.Cells.Item(40, 4).Value = '7.1.0. Moscow'
$ws.Cells.Item(41, 4).Value = '7.2.0. St.Petersburg'
$ws.Cells.Item(42, 4).Value = '7.3.0. 1mln+ population'
$ws.Cells.Item(43, 4).Value = '7.4.0. 500 thd - 1 mln population'
$ws.Cells.Item(44, 4).Value = '7.5.0. 10 thd - 500 thd population'
$ws.Cells.Item(45, 4).Value = '7.6.0. Rural Area'
$ws.Cells.Item(47, 4).Value = '8.1.0. <10000 RBLS'
$ws.Cells.Item(48, 4).Value = '8.2.0. 10001-45000 RBLS'
$ws.Cells.Item(49, 4).Value = '8.3.0. 45001-60000 RBLS'
$ws.Cells.Item(50, 4).Value = '8.4.0. > 60000 RBLS'
$ws.Cells.Item(52, 4).Value = '9.1.0. 0-12 month'
$ws.Cells.Item(53, 4).Value = '9.2.0. 13-24 month'
$ws.Cells.Item(54, 4).Value = '9.3.0. 25-36 month'
$ws.Cells.Item(56, 4).Value = '8.1.0. No teen in HH'
$ws.Cells.Item(57, 4).Value = '8.2.0. Teen girls 10 to 17 y.o'
$ws.Cells.Item(58, 4).Value = '8.3.0. Teen girls 18 to 24 y.o'
$ws.Cells.Item(60, 4).Value = '9.1.0. 18-24 yrs old'
$ws.Cells.Item(61, 4).Value = '9.2.0. 25-34 yrs old'
$ws.Cells.Item(62, 4).Value = '9.3.0. 35-49 yrs old'
$ws.Cells.Item(63, 4).Value = '9.4.0. 50-64 yrs old'
$ws.Cells.Item(64, 4).Value = '9.5.0. 65+ yrs old'

# --- section header ("N.0. Category") rows (in row order) ---
$ws.Cells.Item(3, 4).Value = '1.0. Age'
$ws.Cells.Item(9, 4).Value = '2.0. Size of households'
$ws.Cells.Item(14, 4).Value = '3.0. Children'
$ws.Cells.Item(19, 4).Value = '4.0. Occupation'
$ws.Cells.Item(23, 4).Value = '5.0. Affluency'
$ws.Cells.Item(28, 4).Value = '6.0. Lifestages'
$ws.Cells.Item(39, 4).Value = '7.0. Geographical Area'
$ws.Cells.Item(46, 4).Value = '8.0. Income'
$ws.Cells.Item(51, 4).Value = '9.0. Age of Baby'
$ws.Cells.Item(55, 4).Value = '8.0. POME'
$ws.Cells.Item(59, 4).Value = '9.0. Age of Total household'

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select()
